# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple worksheets, per the scheduled-runner data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1698.5
$ws.Range("I33").Value = 438.1
$ws.Range("K33").Value = 438.1
$ws.Range("M33").Value = -209.1

$ws.Range("H40").Value = 1725.8649
$ws.Range("I40").Value = 1548.5714
$ws.Range("K40").Value = 1548.5714
$ws.Range("M40").Value = -1373.5714

$ws.Range("H74").Value = 4420.8887
$ws.Range("I74").Value = 4336
$ws.Range("J74").Value = 5100
$ws.Range("K74").Value = 4336
$ws.Range("L74").Value = 5100
$ws.Range("M74").Value = -3400
$ws.Range("N74").Value = -6972

$ws.Range("H77").Value = 4420.8887
$ws.Range("I77").Value = 4336
$ws.Range("J77").Value = 5100
$ws.Range("K77").Value = 21680
$ws.Range("L77").Value = 25500
$ws.Range("M77").Value = -17000
$ws.Range("N77").Value = -34860

$ws.Range("H125").Value = 2999.7693
$ws.Range("J125").Value = 2751.875
$ws.Range("L125").Value = 24766.875
$ws.Range("N125").Value = -29686.875


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 61000
$ws.Range("I20").Value = 73678.57000000001
$ws.Range("J20").Value = 1833.3334
$ws.Range("K20").Value = 73678.57000000001
$ws.Range("L20").Value = 1833.3334
$ws.Range("M20").Value = -73431.57000000001
$ws.Range("N20").Value = -2327.3334

$ws.Range("H80").Value = 1565.7778
$ws.Range("J80").Value = 1888.5834
$ws.Range("L80").Value = 1888.5834
$ws.Range("N80").Value = -3884.5834

$ws.Range("H83").Value = 1565.7778
$ws.Range("J83").Value = 1888.5834
$ws.Range("L83").Value = 9442.916999999999
$ws.Range("N83").Value = -19426.917

$ws.Range("H105").Value = 135176.6
$ws.Range("I105").Value = 144535.72
$ws.Range("J105").Value = 126987.375
$ws.Range("K105").Value = 144535.72
$ws.Range("L105").Value = 126987.375
$ws.Range("M105").Value = -142788.72
$ws.Range("N105").Value = -130481.375

$ws.Range("H134").Value = 1839.6888
$ws.Range("I134").Value = 1936.25
$ws.Range("J134").Value = 1453.4445
$ws.Range("K134").Value = 5808.75
$ws.Range("L134").Value = 4360.333500000001
$ws.Range("M134").Value = -3273.75
$ws.Range("N134").Value = -9430.333500000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28651.184
$ws.Range("I31").Value = 43176.875
$ws.Range("K31").Value = 43176.875
$ws.Range("M31").Value = -42881.875

$ws.Range("H34").Value = 28651.184
$ws.Range("I34").Value = 43176.875
$ws.Range("K34").Value = 43176.875
$ws.Range("M34").Value = -42974.875

$ws.Range("H86").Value = 3966.5
$ws.Range("J86").Value = 4159.8
$ws.Range("L86").Value = 4159.8
$ws.Range("N86").Value = -6405.8

$ws.Range("H89").Value = 3966.5
$ws.Range("J89").Value = 4159.8
$ws.Range("L89").Value = 20799
$ws.Range("N89").Value = -32031

$ws.Range("H107").Value = 872.7857
$ws.Range("I107").Value = 1158.7142
$ws.Range("J107").Value = 586.8570999999999
$ws.Range("K107").Value = 1158.7142
$ws.Range("L107").Value = 586.8570999999999
$ws.Range("M107").Value = 761.2858000000001
$ws.Range("N107").Value = -4426.8571

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H122").Value = 694.3333
$ws.Range("I122").Value = 592.7143
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 1778.1429
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = 671.8571000000002
$ws.Range("N122").Value = -8050


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H118").Value = 800
$ws.Range("I118").Value = 800
$ws.Range("K118").Value = 2400
$ws.Range("M118").Value = -1157

$ws.Range("H131").Value = 667703.6
$ws.Range("I131").Value = 559.9231
$ws.Range("J131").Value = 767391.8
$ws.Range("K131").Value = 1679.7693
$ws.Range("L131").Value = 2302175.4
$ws.Range("M131").Value = 3360.2307
$ws.Range("N131").Value = -2312255.4

$ws.Range("H132").Value = 2830.55
$ws.Range("I132").Value = 2964.2727
$ws.Range("J132").Value = 2667.111
$ws.Range("K132").Value = 26678.4543
$ws.Range("L132").Value = 24003.999
$ws.Range("M132").Value = -24148.4543
$ws.Range("N132").Value = -29063.999

$ws.Range("H139").Value = 2472.5
$ws.Range("I139").Value = 1431.1111
$ws.Range("J139").Value = 3097.3333
$ws.Range("K139").Value = 4293.3333
$ws.Range("L139").Value = 9291.999899999999
$ws.Range("M139").Value = 846.6666999999998
$ws.Range("N139").Value = -19571.9999


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 49397.223
$ws.Range("I70").Value = 78646.74000000001
$ws.Range("K70").Value = 78646.74000000001
$ws.Range("M70").Value = -78376.74000000001

$ws.Range("H73").Value = 49397.223
$ws.Range("I73").Value = 78646.74000000001
$ws.Range("K73").Value = 78646.74000000001
$ws.Range("M73").Value = -77710.74000000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1047.1428
$ws.Range("I22").Value = 986
$ws.Range("J22").Value = 1063.8182
$ws.Range("K22").Value = 986
$ws.Range("L22").Value = 1063.8182
$ws.Range("M22").Value = -691
$ws.Range("N22").Value = -1653.8182

$ws.Range("H27").Value = 1047.1428
$ws.Range("I27").Value = 986
$ws.Range("J27").Value = 1063.8182
$ws.Range("K27").Value = 986
$ws.Range("L27").Value = 1063.8182
$ws.Range("M27").Value = -879
$ws.Range("N27").Value = -1277.8182


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 22820
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 22820
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 22820
$ws.Range("N4").Value = -23046
$ws.Range("M4").ClearContents()

$ws.Range("H100").Value = 251400.5
$ws.Range("I100").Value = 334867.34
$ws.Range("K100").Value = 669734.6800000001
$ws.Range("M100").Value = -669193.6800000001

$ws.Range("H107").Value = 200198.6
$ws.Range("I107").Value = 264.66666
$ws.Range("J107").Value = 500099.5
$ws.Range("K107").Value = 793.9999799999999
$ws.Range("L107").Value = 1500298.5
$ws.Range("M107").Value = 1126.00002
$ws.Range("N107").Value = -1504138.5

$ws.Range("H126").Value = 2380.25
$ws.Range("I126").Value = 2491.1333
$ws.Range("J126").Value = 2047.6
$ws.Range("K126").Value = 7473.3999
$ws.Range("L126").Value = 6142.799999999999
$ws.Range("M126").Value = -5003.3999
$ws.Range("N126").Value = -11082.8

